# Comparison table fix (slide "Comparison table"):
# The "Progress for Readers" row, "Double Instance Locking" column was
# incorrectly marked "Lock-Free" (green) -- it should say "Blocking" (red),
# matching the other "Blocking" cells in the table.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

$cell = $tbl.Cell(3, 4)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "Blocking"
$tr.Font.Color.RGB = 0x0000FF
